$wb = $excel.ActiveWorkbook

# Sheet ALC, row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value2 = 9824.888999999999
$ws.Range("I103").Value2 = 428
$ws.Range("J103").Value2 = 85000
$ws.Range("K103").Value2 = 1284
$ws.Range("L103").Value2 = 255000
$ws.Range("M103").Value2 = -698
$ws.Range("N103").Value2 = -256172

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value2 = 1402628.8
$ws.Range("I137").Value2 = 1643456.1
$ws.Range("J137").Value2 = 5830.6
$ws.Range("K137").Value2 = 4930368.300000001
$ws.Range("L137").Value2 = 17491.8
$ws.Range("M137").Value2 = -4927818.300000001
$ws.Range("N137").Value2 = -22591.8

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value2 = 2908.6155
$ws.Range("I138").Value2 = 2106.5
$ws.Range("K138").Value2 = 6319.5
$ws.Range("M138").Value2 = -1179.5

# Sheet ARM, row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value2 = 2519
$ws.Range("I45").Value2 = 2778.6667
$ws.Range("J45").Value2 = 1740
$ws.Range("K45").Value2 = 2778.6667
$ws.Range("L45").Value2 = 1740
$ws.Range("M45").Value2 = -2401.6667
$ws.Range("N45").Value2 = -2494

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 2349.4285
$ws.Range("I61").Value2 = 1907.6666
$ws.Range("J61").Value2 = 5000
$ws.Range("K61").Value2 = 1907.6666
$ws.Range("L61").Value2 = 5000
$ws.Range("M61").Value2 = -1695.6666
$ws.Range("N61").Value2 = -5424

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value2 = 3906.4375
$ws.Range("I74").Value2 = 3441.9167
$ws.Range("J74").Value2 = 5300
$ws.Range("K74").Value2 = 3441.9167
$ws.Range("L74").Value2 = 5300
$ws.Range("M74").Value2 = -2567.9167
$ws.Range("N74").Value2 = -7048

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value2 = 3906.4375
$ws.Range("I77").Value2 = 3441.9167
$ws.Range("J77").Value2 = 5300
$ws.Range("K77").Value2 = 17209.5835
$ws.Range("L77").Value2 = 26500
$ws.Range("M77").Value2 = -12841.5835
$ws.Range("N77").Value2 = -35236

# Sheet ARM, row 109
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value2 = 33204.76
$ws.Range("J109").Value2 = 33204.76
$ws.Range("L109").Value2 = 33204.76
$ws.Range("N109").Value2 = -35978.76

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value2 = 2349.4285
$ws.Range("I136").Value2 = 1907.6666
$ws.Range("J136").Value2 = 5000
$ws.Range("K136").Value2 = 5722.9998
$ws.Range("L136").Value2 = 15000
$ws.Range("M136").Value2 = -3172.9998
$ws.Range("N136").Value2 = -20100

# Sheet BSM, row 132
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value2 = 56056.832
$ws.Range("J132").Value2 = 56056.832
$ws.Range("L132").Value2 = 56056.832
$ws.Range("N132").Value2 = -66176.83199999999

# Sheet CRP, row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value2 = 3786.5715
$ws.Range("I86").Value2 = 3569
$ws.Range("K86").Value2 = 3569
$ws.Range("M86").Value2 = -2446

# Sheet CRP, row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value2 = 3786.5715
$ws.Range("I89").Value2 = 3569
$ws.Range("K89").Value2 = 17845
$ws.Range("M89").Value2 = -12229

# Sheet CRP, row 137
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H137").Value2 = 45353.332
$ws.Range("J137").Value2 = 45353.332
$ws.Range("L137").Value2 = 45353.332
$ws.Range("N137").Value2 = -55553.332

# Sheet CUL, row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 495699.88
$ws.Range("J5").Value2 = 787019.7
$ws.Range("L5").Value2 = 2361059.1
$ws.Range("N5").Value2 = -2361283.1

# Sheet CUL, row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value2 = 10195.448
$ws.Range("I34").Value2 = 20000.334
$ws.Range("J34").Value2 = 7637.6523
$ws.Range("K34").Value2 = 60001.00199999999
$ws.Range("L34").Value2 = 22912.9569
$ws.Range("M34").Value2 = -59917.00199999999
$ws.Range("N34").Value2 = -23080.9569

# Sheet CUL, row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value2 = 10375.2
$ws.Range("J39").Value2 = 10628.41
$ws.Range("L39").Value2 = 31885.23
$ws.Range("N39").Value2 = -32473.23

# Sheet CUL, row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value2 = 4903.5713
$ws.Range("I55").Value2 = 1000
$ws.Range("J55").Value2 = 5203.846
$ws.Range("K55").Value2 = 3000
$ws.Range("L55").Value2 = 15611.538
$ws.Range("M55").Value2 = -2823
$ws.Range("N55").Value2 = -15965.538

# Sheet CUL, row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value2 = 6579547
$ws.Range("J113").Value2 = 13889442
$ws.Range("L113").Value2 = 41668326
$ws.Range("N113").Value2 = -41672666

# Sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value2 = 807.35
$ws.Range("I131").Value2 = 308
$ws.Range("J131").Value2 = 833.6316
$ws.Range("K131").Value2 = 924
$ws.Range("L131").Value2 = 2500.8948
$ws.Range("M131").Value2 = 4116
$ws.Range("N131").Value2 = -12580.8948

# Sheet CUL, row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value2 = 495699.88
$ws.Range("J135").Value2 = 787019.7
$ws.Range("L135").Value2 = 7083177.3
$ws.Range("N135").Value2 = -7088247.3

# Sheet GSM, row 11
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value2 = 10468825
$ws.Range("I11").Value2 = 25000028
$ws.Range("J11").Value2 = 781355.7
$ws.Range("K11").Value2 = 25000028
$ws.Range("L11").Value2 = 781355.7
$ws.Range("M11").Value2 = -24999889
$ws.Range("N11").Value2 = -781633.7

# Sheet GSM, row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value2 = 3943.6667
$ws.Range("I102").Value2 = 2734.7
$ws.Range("J102").Value2 = 9988.5
$ws.Range("K102").Value2 = 2734.7
$ws.Range("L102").Value2 = 9988.5
$ws.Range("M102").Value2 = -1112.7
$ws.Range("N102").Value2 = -13232.5

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value2 = 2843.8965
$ws.Range("I132").Value2 = 2044.7916
$ws.Range("J132").Value2 = 6679.6
$ws.Range("K132").Value2 = 6134.3748
$ws.Range("L132").Value2 = 20038.8
$ws.Range("M132").Value2 = -3604.3748
$ws.Range("N132").Value2 = -25098.8

# Sheet LTW, row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 1187.1111
$ws.Range("I16").Value2 = 1187.1111
$ws.Range("J16").Value2 = 0
$ws.Range("K16").Value2 = 1187.1111
$ws.Range("L16").Value2 = 0
$ws.Range("M16").Value2 = -1017.1111
$ws.Range("N16").ClearContents()

# Sheet LTW, row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value2 = 2135.9375
$ws.Range("I46").Value2 = 1802.7778
$ws.Range("J46").Value2 = 2564.2856
$ws.Range("K46").Value2 = 1802.7778
$ws.Range("L46").Value2 = 2564.2856
$ws.Range("M46").Value2 = -1614.7778
$ws.Range("N46").Value2 = -2940.2856

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value2 = 9750
$ws.Range("I122").Value2 = 8875
$ws.Range("J122").Value2 = 11500
$ws.Range("K122").Value2 = 26625
$ws.Range("L122").Value2 = 34500
$ws.Range("M122").Value2 = -24175
$ws.Range("N122").Value2 = -39400

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value2 = 5647.516
$ws.Range("I122").Value2 = 3958
$ws.Range("J122").Value2 = 9195.5
$ws.Range("K122").Value2 = 11874
$ws.Range("L122").Value2 = 27586.5
$ws.Range("M122").Value2 = -9424
$ws.Range("N122").Value2 = -32486.5

# Sheet WVR, row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value2 = 2037.5667
$ws.Range("I126").Value2 = 1121.0476
$ws.Range("J126").Value2 = 4176.1113
$ws.Range("K126").Value2 = 3363.142800000001
$ws.Range("L126").Value2 = 12528.3339
$ws.Range("M126").Value2 = -893.1428000000005
$ws.Range("N126").Value2 = -17468.3339
